$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7676553333333332
$ws.Range("H2").Value = 2.302966
$ws.Range("I2").Value = 0.3736977786965754
$ws.Range("J2").Value = 0.3736977786965754
$ws.Range("M2").Value = 43.97948166666666
$ws.Range("N2").Value = 131.938445
$ws.Range("O2").Value = 0.3260725128076164
$ws.Range("P2").Value = 0.3260725128076164
$ws.Range("Q2").Value = 33.76108365865222
$ws.Range("R2").Value = 303.8497529278699
$ws.Range("S2").Value = 0.1218525737302169
$ws.Range("T2").Value = 0.1218525737302169

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7676553333333332
$ws.Range("H3").Value = 2.302966
$ws.Range("I3").Value = 0.3736977786965754
$ws.Range("J3").Value = 0.3736977786965754
$ws.Range("M3").Value = 57.80064033333333
$ws.Range("O3").Value = 0.4285452970598356
$ws.Range("P3").Value = 0.4285452970598356
$ws.Range("Q3").Value = 44.3709698219651
$ws.Range("R3").Value = 399.3387283976859
$ws.Range("S3").Value = 0.1601464255821246
$ws.Range("T3").Value = 0.1601464255821246

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7676553333333332
$ws.Range("H4").Value = 2.302966
$ws.Range("I4").Value = 0.3736977786965754
$ws.Range("J4").Value = 0.3736977786965754
$ws.Range("M4").Value = 20.92900166666667
$ws.Range("N4").Value = 62.787005
$ws.Range("O4").Value = 0.15517172793733
$ws.Range("P4").Value = 0.15517172793733
$ws.Range("Q4").Value = 16.06625975075889
$ws.Range("R4").Value = 144.59633775683
$ws.Range("S4").Value = 0.05798733004668956
$ws.Range("T4").Value = 0.05798733004668954

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7676553333333332
$ws.Range("H5").Value = 2.302966
$ws.Range("I5").Value = 0.3736977786965754
$ws.Range("J5").Value = 0.3736977786965754
$ws.Range("M5").Value = 12.167261
$ws.Range("N5").Value = 36.501783
$ws.Range("O5").Value = 0.0902104621952179
$ws.Range("P5").Value = 0.0902104621952179
$ws.Range("Q5").Value = 9.340262798708666
$ws.Range("R5").Value = 84.062365188378
$ws.Range("S5").Value = 0.03371144933754432
$ws.Range("T5").Value = 0.03371144933754432

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.811404
$ws.Range("H6").Value = 2.434212
$ws.Range("I6").Value = 0.3949948098567449
$ws.Range("J6").Value = 0.3949948098567449
$ws.Range("M6").Value = 43.97948166666666
$ws.Range("N6").Value = 131.938445
$ws.Range("O6").Value = 0.3260725128076164
$ws.Range("P6").Value = 0.3260725128076164
$ws.Range("Q6").Value = 35.68512734226
$ws.Range("R6").Value = 321.16614608034
$ws.Range("S6").Value = 0.1287969501959555
$ws.Range("T6").Value = 0.1287969501959555

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.811404
$ws.Range("H7").Value = 2.434212
$ws.Range("I7").Value = 0.3949948098567449
$ws.Range("J7").Value = 0.3949948098567449
$ws.Range("M7").Value = 57.80064033333333
$ws.Range("O7").Value = 0.4285452970598356
$ws.Range("P7").Value = 0.4285452970598356
$ws.Range("Q7").Value = 46.89967076902799
$ws.Range("R7").Value = 422.097036921252
$ws.Range("S7").Value = 0.169273168127152
$ws.Range("T7").Value = 0.169273168127152

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.811404
$ws.Range("H8").Value = 2.434212
$ws.Range("I8").Value = 0.3949948098567449
$ws.Range("J8").Value = 0.3949948098567449
$ws.Range("M8").Value = 20.92900166666667
$ws.Range("N8").Value = 62.787005
$ws.Range("O8").Value = 0.15517172793733
$ws.Range("P8").Value = 0.15517172793733
$ws.Range("Q8").Value = 16.98187566834
$ws.Range("R8").Value = 152.83688101506
$ws.Range("S8").Value = 0.06129202717174821
$ws.Range("T8").Value = 0.0612920271717482

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.811404
$ws.Range("H9").Value = 2.434212
$ws.Range("I9").Value = 0.3949948098567449
$ws.Range("J9").Value = 0.3949948098567449
$ws.Range("M9").Value = 12.167261
$ws.Range("N9").Value = 36.501783
$ws.Range("O9").Value = 0.0902104621952179
$ws.Range("P9").Value = 0.0902104621952179
$ws.Range("Q9").Value = 9.872564244444002
$ws.Range("R9").Value = 88.85307819999601
$ws.Range("S9").Value = 0.03563266436188917
$ws.Range("T9").Value = 0.03563266436188917

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.475155
$ws.Range("H10").Value = 1.425465
$ws.Range("I10").Value = 0.2313074114466796
$ws.Range("J10").Value = 0.2313074114466796
$ws.Range("M10").Value = 43.97948166666666
$ws.Range("N10").Value = 131.938445
$ws.Range("O10").Value = 0.3260725128076164
$ws.Range("P10").Value = 0.3260725128076164
$ws.Range("Q10").Value = 20.897070611325
$ws.Range("R10").Value = 188.073635501925
$ws.Range("S10").Value = 0.07542298888144403
$ws.Range("T10").Value = 0.07542298888144403

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.475155
$ws.Range("H11").Value = 1.425465
$ws.Range("I11").Value = 0.2313074114466796
$ws.Range("J11").Value = 0.2313074114466796
$ws.Range("M11").Value = 57.80064033333333
$ws.Range("O11").Value = 0.4285452970598356
$ws.Range("P11").Value = 0.4285452970598356
$ws.Range("Q11").Value = 27.464263257585
$ws.Range("R11").Value = 247.178369318265
$ws.Range("S11").Value = 0.09912570335055894
$ws.Range("T11").Value = 0.09912570335055894

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.475155
$ws.Range("H12").Value = 1.425465
$ws.Range("I12").Value = 0.2313074114466796
$ws.Range("J12").Value = 0.2313074114466796
$ws.Range("M12").Value = 20.92900166666667
$ws.Range("N12").Value = 62.787005
$ws.Range("O12").Value = 0.15517172793733
$ws.Range("P12").Value = 0.15517172793733
$ws.Range("Q12").Value = 9.944519786925001
$ws.Range("R12").Value = 89.500678082325
$ws.Range("S12").Value = 0.03589237071889222
$ws.Range("T12").Value = 0.03589237071889221

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.475155
$ws.Range("H13").Value = 1.425465
$ws.Range("I13").Value = 0.2313074114466796
$ws.Range("J13").Value = 0.2313074114466796
$ws.Range("M13").Value = 12.167261
$ws.Range("N13").Value = 36.501783
$ws.Range("O13").Value = 0.0902104621952179
$ws.Range("P13").Value = 0.0902104621952179
$ws.Range("Q13").Value = 5.781334900455001
$ws.Range("R13").Value = 52.032014104095
$ws.Range("S13").Value = 0.0208663484957844
$ws.Range("T13").Value = 0.0208663484957844
